$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Insert a new worksheet "2022-Q1" right before the "总计" sheet.
# ---------------------------------------------------------------------------
$templateSheet = $wb.Worksheets.Item("2021-Q4")
$totalSheetBefore = $wb.Worksheets.Item("总计")

$newSheet = $wb.Worksheets.Add($totalSheetBefore)
$newSheet.Name = "2022-Q1"

# NOTE: in this runtime, the object returned by Worksheets.Add(Before) and
# the "Before" reference itself end up pointing at the *same* (newly
# inserted) sheet, so the original "总计" sheet must be looked up again by
# name afterwards to get a handle on the real summary sheet.
$totalSheet = $wb.Worksheets.Item("总计")

# Copy the header formatting (bold, centered, bordered) from an existing
# quarter sheet so the new sheet matches the established look.
$templateSheet.Range("B1:H1").Copy()
$newSheet.Range("B1:H1").PasteSpecial(-4122)
$newSheet.Application.CutCopyMode = $false

$newSheet.Range("B1").Value = "基金代码"
$newSheet.Range("C1").Value = "基金名称"
$newSheet.Range("D1").Value = "基金规模"
$newSheet.Range("E1").Value = "股票总仓位"
$newSheet.Range("F1").Value = "仓位占比"
$newSheet.Range("G1").Value = "持有市值(亿元)"
$newSheet.Range("H1").Value = "仓位排名"

# Row 2 / Row 3 index column (A) uses the same bold/border style as the
# header and holds a plain running number (0, 1, ...).
$templateSheet.Range("A2").Copy()
$newSheet.Range("A2:A3").PasteSpecial(-4122)
$newSheet.Application.CutCopyMode = $false

$newSheet.Range("A2").Value = 0
$newSheet.Range("A3").Value = 1

# Columns B:G hold numeric-looking text (fund codes / percentages), so
# force a text number format before assigning the values, otherwise Excel
# would silently convert them into real numbers and lose leading zeros.
$textRange = $newSheet.Range("B2:G3")
$textRange.NumberFormat = "@"

$newSheet.Range("B2").Value = "005051"
$newSheet.Range("C2").Value = "上投摩根标普港股通低波红利指数A"
$newSheet.Range("D2").Value = "4.02"
$newSheet.Range("E2").Value = "92.23"
$newSheet.Range("F2").Value = "4.26"
$newSheet.Range("G2").Value = "0.1713"

$newSheet.Range("B3").Value = "005052"
$newSheet.Range("C3").Value = "上投摩根标普港股通低波红利指数C"
$newSheet.Range("D3").Value = "2.61"
$newSheet.Range("E3").Value = "92.23"
$newSheet.Range("F3").Value = "4.26"
$newSheet.Range("G3").Value = "0.1112"

$textRange.Style = "Normal"

$newSheet.Range("H2").Value = 2
$newSheet.Range("H3").Value = 2

# ---------------------------------------------------------------------------
# 2. Add a new "2022-Q1" row on top of the "总计" (summary) sheet data and
#    keep the running index (column A) and the existing rows in sync.
# ---------------------------------------------------------------------------
$totalSheet.Rows.Item(2).Insert()

# The row that used to be row 2 is now row 3; reuse its (known-good)
# formatting for the brand-new row 2 instead of whatever Insert() guessed.
$totalSheet.Range("A3").Copy()
$totalSheet.Range("A2").PasteSpecial(-4122)
$totalSheet.Application.CutCopyMode = $false

$totalSheet.Range("B3:D3").Copy()
$totalSheet.Range("B2:D2").PasteSpecial(-4122)
$totalSheet.Application.CutCopyMode = $false

$totalSheet.Range("A2").Value = 0
$totalSheet.Range("B2").Value = "2022-Q1"
$totalSheet.Range("C2").Value = 2
$totalSheet.Range("D2").Value = 0.28

$totalSheet.Range("A3").Value = 1
$totalSheet.Range("A4").Value = 2
$totalSheet.Range("A5").Value = 3
$totalSheet.Range("A6").Value = 4
$totalSheet.Range("A7").Value = 5
